$d = $word.ActiveDocument

$d.Content.Find.Execute("43÷4=10, 3", $true, $false, $false, $false, $false, $true, 1, $false, "61÷2=30, 1", 2)
$d.Content.Find.Execute("32÷6=5, 2", $true, $false, $false, $false, $false, $true, 1, $false, "41÷4=10, 1", 2)
$d.Content.Find.Execute("93÷7=13, 2", $true, $false, $false, $false, $false, $true, 1, $false, "86÷2=43, 0", 2)
$d.Content.Find.Execute("26÷3=8, 2", $true, $false, $false, $false, $false, $true, 1, $false, "44÷5=8, 4", 2)
$d.Content.Find.Execute("62÷9=6, 8", $true, $false, $false, $false, $false, $true, 1, $false, "15÷5=3, 0", 2)
$d.Content.Find.Execute("71÷2=35, 1", $true, $false, $false, $false, $false, $true, 1, $false, "98÷3=32, 2", 2)
$d.Content.Find.Execute("91÷9=10, 1", $true, $false, $false, $false, $false, $true, 1, $false, "94÷7=13, 3", 2)
$d.Content.Find.Execute("45÷3=15, 0", $true, $false, $false, $false, $false, $true, 1, $false, "27÷6=4, 3", 2)
$d.Content.Find.Execute("48÷9=5, 3", $true, $false, $false, $false, $false, $true, 1, $false, "77÷2=38, 1", 2)
$d.Content.Find.Execute("45÷7=6, 3", $true, $false, $false, $false, $false, $true, 1, $false, "70÷3=23, 1", 2)
$d.Content.Find.Execute("48÷5=9, 3", $true, $false, $false, $false, $false, $true, 1, $false, "91÷4=22, 3", 2)
$d.Content.Find.Execute("74÷6=12, 2", $true, $false, $false, $false, $false, $true, 1, $false, "59÷4=14, 3", 2)
$d.Content.Find.Execute("77÷9=8, 5", $true, $false, $false, $false, $false, $true, 1, $false, "65÷2=32, 1", 2)
$d.Content.Find.Execute("47÷6=7, 5", $true, $false, $false, $false, $false, $true, 1, $false, "74÷9=8, 2", 2)
$d.Content.Find.Execute("63÷5=12, 3", $true, $false, $false, $false, $false, $true, 1, $false, "29÷5=5, 4", 2)
$d.Content.Find.Execute("55÷7=7, 6", $true, $false, $false, $false, $false, $true, 1, $false, "89÷3=29, 2", 2)
$d.Content.Find.Execute("35÷9=3, 8", $true, $false, $false, $false, $false, $true, 1, $false, "12÷6=2, 0", 2)
$d.Content.Find.Execute("64÷3=21, 1", $true, $false, $false, $false, $false, $true, 1, $false, "10÷7=1, 3", 2)
$d.Content.Find.Execute("82÷4=20, 2", $true, $false, $false, $false, $false, $true, 1, $false, "47÷9=5, 2", 2)
$d.Content.Find.Execute("30÷5=6, 0", $true, $false, $false, $false, $false, $true, 1, $false, "59÷5=11, 4", 2)
$d.Content.Find.Execute("79÷8=9, 7", $true, $false, $false, $false, $false, $true, 1, $false, "11÷3=3, 2", 2)
$d.Content.Find.Execute("94÷2=47, 0", $true, $false, $false, $false, $false, $true, 1, $false, "60÷8=7, 4", 2)
$d.Content.Find.Execute("23÷9=2, 5", $true, $false, $false, $false, $false, $true, 1, $false, "94÷8=11, 6", 2)
$d.Content.Find.Execute("46÷6=7, 4", $true, $false, $false, $false, $false, $true, 1, $false, "58÷9=6, 4", 2)
$d.Content.Find.Execute("42÷8=5, 2", $true, $false, $false, $false, $false, $true, 1, $false, "65÷6=10, 5", 2)
